$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 26 ("RM 232") entirely - remaining rows shift up by one.
$ws.Rows.Item(26).Delete()

# After that deletion, the former row 28 ("SC 92") is now row 27; delete it too.
$ws.Rows.Item(27).Delete()

# A couple of the "Missing at random" cells in column E were (re)imputed with
# different values as part of this pass - fix those up on the now-shifted rows.
$ws.Range("E26").Value = ""
$ws.Range("E27").Value = -10
$ws.Range("E29").Value = ""

# The used range now only goes down to row 33.
$ws.Range("A1:F33").Select()
